# Update TPM-derived values in the NATMI LR-pairs output sheet.
# Columns M:T for rows 2-5 are replaced with newly recomputed values
# (receptor/edge expression values & derived specificity scores),
# per commit "update scripts wuth new tpm".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("M2").Value = 0.8908616666666668
$ws.Range("N2").Value = 2.672585
$ws.Range("O2").Value = 0.04079002072021364
$ws.Range("P2").Value = 0.04079002072021363
$ws.Range("Q2").Value = 0.4254874101555556
$ws.Range("R2").Value = 3.829386691400001
$ws.Range("S2").Value = 0.04079002072021364
$ws.Range("T2").Value = 0.04079002072021363

# Row 3
$ws.Range("M3").Value = 16.81477433333333
$ws.Range("N3").Value = 50.444323
$ws.Range("O3").Value = 0.7699006693471485
$ws.Range("P3").Value = 0.7699006693471484
$ws.Range("Q3").Value = 8.030960418591112
$ws.Range("R3").Value = 72.27864376732001
$ws.Range("S3").Value = 0.7699006693471485
$ws.Range("T3").Value = 0.7699006693471484

# Row 4
$ws.Range("M4").Value = 3.879966
$ws.Range("N4").Value = 11.639898
$ws.Range("O4").Value = 0.1776526024808091
$ws.Range("P4").Value = 0.1776526024808091
$ws.Range("Q4").Value = 1.85312349448
$ws.Range("R4").Value = 16.67811145032
$ws.Range("S4").Value = 0.1776526024808091
$ws.Range("T4").Value = 0.1776526024808091

# Row 5
$ws.Range("M5").Value = 0.2545846666666667
$ws.Range("N5").Value = 0.763754
$ws.Range("O5").Value = 0.01165670745182886
$ws.Range("P5").Value = 0.01165670745182886
$ws.Range("Q5").Value = 0.1215930312622222
$ws.Range("R5").Value = 1.09433728136
$ws.Range("S5").Value = 0.01165670745182886
$ws.Range("T5").Value = 0.01165670745182886
